# This edit corresponds to the workbook being re-saved by Excel after the
# project's dataset-spec generator was switched from a Jupyter notebook to
# a plain Python script (run via "make specs"). The sheet data, layout and
# styles are unchanged - the regenerated file simply picks up Excel's
# default page-margin metadata instead of the previous custom values.
#
# Previous (custom) margins -> Excel defaults:
#   left/right : 0.75in -> 0.7in
#   top/bottom : 1in    -> 0.75in
#   header/footer : 0.5in -> 0.3in

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.PageSetup.LeftMargin   = 50.4   # 0.7in in points (72 pts/in)
$ws.PageSetup.RightMargin  = 50.4
$ws.PageSetup.TopMargin    = 54.0   # 0.75in
$ws.PageSetup.BottomMargin = 54.0
$ws.PageSetup.HeaderMargin = 21.6   # 0.3in
$ws.PageSetup.FooterMargin = 21.6
